$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 7791.630278080007
$ws.Range("E2").Value = 231859.0077278444
$ws.Range("I2").Value = 129402.25083824
$ws.Range("L2").Value = 388196.236028304
$ws.Range("M2").Value = 87070.39609781333
$ws.Range("N2").Value = 56665.56446366481
$ws.Range("O2").Value = 55688.35379128032

$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 37908.85114466308
$ws.Range("E2").Value = 216989.4393747085
$ws.Range("I2").Value = 224340.9369390888
$ws.Range("L2").Value = 147536.3839604098
$ws.Range("M2").Value = 93673.10875507205
$ws.Range("N2").Value = 27145.4596892804
$ws.Range("O2").Value = 40388.37785974222

$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 22895.69120990697
$ws.Range("B2").Value = 18514.63588195528
$ws.Range("E2").Value = 89533.4725380327
$ws.Range("I2").Value = 120308.2182965601
$ws.Range("M2").Value = 28613.58482505595
$ws.Range("N2").Value = 35950.49127023212
$ws.Range("O2").Value = 21550.65044883481

$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 811.8131926239505

$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 27527.29374076518
$ws.Range("N2").Value = 4146.158925555366
$ws.Range("O2").Value = 18378.03620052724
